$d = $word.ActiveDocument
$bmNames = New-Object System.Collections.ArrayList
$bmCounter = 0

function Add-Text {
    param($para, [string]$txt)
    $r = $para.Range
    $r.Collapse(0)
    $r.InsertAfter($txt)
}

function Add-Split {
    param($para)
    $script:bmCounter += 1
    $name = "zzsplit" + $script:bmCounter
    $mid = $para.Range
    $mid.Collapse(0)
    $d.Bookmarks.Add($name, $mid) | Out-Null
    [void]$bmNames.Add($name)
}

function Add-LineBreak {
    param($para)
    $r = $para.Range
    $r.Collapse(0)
    $r.InsertBreak(6) | Out-Null
}

function New-ListParagraph {
    param($afterPara, [int]$ilvl)
    $r = $afterPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Range.ListFormat.ListLevelNumber = $ilvl + 1
    return $p
}

$cur = $d.Paragraphs.Last

# --- paragraph 0 (ilvl=0) ---
$cur = New-ListParagraph $cur 0
Add-Text $cur 'AMMEND COMMITS'

# --- paragraph 1 (ilvl=1) ---
$cur = New-ListParagraph $cur 1
Add-Text $cur 'Suppos'
Add-Split $cur
Add-Text $cur 'e you just made a commit and then realized you forgot to include a file! '
Add-Split $cur
Add-Text $cur 'Or,'
Add-Split $cur
Add-Text $cur ' maybe you made a typo in the commit message that you want to correct. '

# --- paragraph 2 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur 'Rather than making a '
Add-Split $cur
Add-Text $cur 'brand new'
Add-Split $cur
Add-Text $cur ' separate commit you can “redo” the previous commit using the –amend option'

# --- paragraph 3 (ilvl=3) ---
$cur = New-ListParagraph $cur 3
Add-Text $cur 'Git commit -m “some commit”'

# --- paragraph 4 (ilvl=3) ---
$cur = New-ListParagraph $cur 3
Add-Text $cur 'Git add '
Add-Split $cur
Add-Text $cur 'forgotten_file'

# --- paragraph 5 (ilvl=3) ---
$cur = New-ListParagraph $cur 3
Add-Text $cur 'Git commit --amend '

# --- paragraph 6 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur 'NOTE – The Git Commit amend only amends the latest commit not the previous commits'

# --- paragraph 7 (ilvl=0) ---
$cur = New-ListParagraph $cur 0
Add-Text $cur 'IGNORING FILES'

# --- paragraph 8 (ilvl=1) ---
$cur = New-ListParagraph $cur 1
Add-Text $cur 'We can tell git which files and directories to ignore '
Add-Split $cur
Add-Text $cur 'in a given'
Add-Split $cur
Add-Text $cur ' '
Add-Split $cur
Add-Text $cur 'repostiroy'
Add-Split $cur
Add-Text $cur ' ,'
Add-Split $cur
Add-Text $cur ' using '
Add-Split $cur
Add-Text $cur 'a .'
Add-Split $cur
Add-Text $cur 'gitignore'
Add-Split $cur
Add-Text $cur ' file. This is useful for files you NEVER want to commit including:'

# --- paragraph 9 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur 'Secrtes'
Add-Split $cur
Add-Text $cur ', API Keys, Credentials, '
Add-Split $cur
Add-Text $cur 'etc'

# --- paragraph 10 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur 'Operationg'
Add-Split $cur
Add-Text $cur ' System Files (.'
Add-Split $cur
Add-Text $cur 'DS_Store'
Add-Split $cur
Add-Text $cur ' on Mac)'

# --- paragraph 11 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur 'Log Files'

# --- paragraph 12 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur 'Dependencies and Packages'

# --- paragraph 13 (ilvl=1) ---
$cur = New-ListParagraph $cur 1
Add-Text $cur '.'
Add-Split $cur
Add-Text $cur 'gitignore'
Add-Split $cur
Add-Text $cur ' – Create a file '
Add-Split $cur
Add-Text $cur 'called .'
Add-Split $cur
Add-Text $cur 'gitignore'
Add-Split $cur
Add-Text $cur ' in the root of a repository. Inside the file, we can write patterns to tell Git which files and folders to ignore:'

# --- paragraph 14 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur '.'
Add-Split $cur
Add-Text $cur 'DS_Store'
Add-Split $cur
Add-Text $cur ' will ignore files name .'
Add-Split $cur
Add-Text $cur 'DS_Store'

# --- paragraph 15 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur 'folderName'
Add-Split $cur
Add-Text $cur ' / will ignore an entire directory'

# --- paragraph 16 (ilvl=2) ---
$cur = New-ListParagraph $cur 2
Add-Text $cur '*.log will ignore any files with the .log extension'

# --- paragraph 17 (ilvl=2) ---
$cur = New-ListParagraph $cur 2

# --- paragraph 18 (ilvl=0) ---
$cur = New-ListParagraph $cur 0
Add-Text $cur 'Gitignore.io file that needs to be added'

# --- paragraph 19 (ilvl=1) ---
$cur = New-ListParagraph $cur 1
Add-Text $cur '.'
Add-Split $cur
Add-Text $cur 'gitignore'
Add-Split $cur
Add-Text $cur ' should be added to the '
Add-Split $cur
Add-Text $cur 'base folder'

# --- paragraph 20 (ilvl=1) ---
$cur = New-ListParagraph $cur 1
Add-Text $cur 'Use .gitignore.io to copy paste general files that needs to be out of io section. '
Add-Split $cur
Add-LineBreak $cur

# clean up all split bookmarks at the very end so runs remain split
foreach ($name in $bmNames) {
    $d.Bookmarks($name).Delete()
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
Write-Host "Remaining bookmarks:" $d.Bookmarks.Count
